# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with the latest scraped values. Several Price values look
# like plain numbers (e.g. "217.23"), so NumberFormat is forced to "@"
# (text) immediately before those assignments to keep them stored as
# text, matching the rest of the column (they'd otherwise be coerced to
# floating point numbers by Excel's automatic type detection).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.248.84'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '1.645.43'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.23'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  +1.88%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('E10').Value = '  +1.81%  '
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = '1.875.43'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '1.634.47'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('E15').Value = '  +3.26%  '
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').Value = '27.231.08'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '219.88'
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.87'
$ws.Range('E21').Value = '  +4.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.57'
$ws.Range('E22').Value = '  +7.25%  '
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.71'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.81'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0511'
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  +1.72%  '
$ws.Range('D35').Value = '1.285.83'
$ws.Range('E35').Value = '  +3.59%  '
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.863'
$ws.Range('E38').Value = '  +4.53%  '
$ws.Range('E39').Value = '  +1.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.807'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').Value = '1.785.63'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.93'
$ws.Range('E45').Value = '  +2.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.97'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('E47').Value = '  +2.07%  '
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.68'
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0973'
$ws.Range('E51').Value = '  +0.57%  '
